$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6190
$ws.Range("C20").Value = 982
$ws.Range("D20").Value = 5584388
$ws.Range("E20").Value = 902.1628432956381
$ws.Range("F20").Value = 6.926930385213326
$ws.Range("G20").Value = 4.24628450106157
$ws.Range("H20").Value = 26.32438544990774
